$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Mod List Database" sheet: insert the new "Conduits Prevent Drowned"
#    mod entry. The sheet is kept sorted (ascending) by Genre (column C), and
#    the new genre text "Tweaks" sorts right after "Tweak" and before
#    "UI/UX", landing the new row at row 65 (pushing everything from the old
#    row 65 down by one).
# ---------------------------------------------------------------------------
$wsDb = $wb.Worksheets.Item("Mod List Database")
$wsDb.Rows("65:65").Insert()

# Set values in the same order the new shared strings were introduced in the
# authored workbook: Mod Name, Description, Genre, Size, Addons.
$wsDb.Range("B65").Value = "Conduits Prevent Drowned"
$wsDb.Range("E65").Value = "Conduits prevent spawning of drowned within range"
$wsDb.Range("C65").Value = "Tweaks"
$wsDb.Range("D65").Value = "Tiny"
$wsDb.Range("F65").Value = "N/A"

# ---------------------------------------------------------------------------
# 2) "Overview" sheet: the curated Tweaks section gets the same new entry,
#    inserted right before the "No hostiles around campfire" row.
# ---------------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")
$wsOv.Rows("147:147").Insert()
$wsOv.Range("B147").Value = "Conduits Prevent Drowned"
$wsOv.Range("C147").Value = "Conduits prevent spawning of drowned within range"

# ---------------------------------------------------------------------------
# 3) View-state bookkeeping to mirror the author's final on-screen state:
#    the "Mod List Database" tab ended up active/selected, scrolled near the
#    newly-inserted row, with the Overview tab scrolled to where its own new
#    row was inserted.
# ---------------------------------------------------------------------------
$wsOv.Select()
$excel.ActiveWindow.ScrollRow = 133
$wsOv.Range("B147:C147").Select()

$wsDb.Select()
$excel.ActiveWindow.ScrollRow = 74
$wsDb.Range("E103").Select()
